$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns touched by the swap: B(2), F(6) through AC(29). Column A (index 1) and C:E (3-5) are left untouched.
$cols = @(2,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)
$rows = @(34,35,88,89,100,101,106,107,114,115,127,128,143,145,149,150,151,162,163,165,166,188,189,214,215,239,240,241)

# 1) Snapshot current (pre-edit) values for every affected row/column so the
#    later writes never read an already-overwritten cell.
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# 2) Row -> row content now sourced from (per commit "Atualizacao de bases das ligas")
$rowSource = @{
    34 = 35
    35 = 34
    88 = 89
    89 = 88
    100 = 101
    101 = 100
    106 = 107
    107 = 106
    114 = 115
    115 = 114
    127 = 128
    128 = 127
    143 = 145
    145 = 143
    149 = 150
    150 = 151
    151 = 149
    162 = 163
    163 = 162
    165 = 166
    166 = 165
    188 = 189
    189 = 188
    214 = 215
    215 = 214
    239 = 240
    240 = 241
    241 = 239
}

# 3) Write each row's cells using the snapshot of its source row
foreach ($r in $rows) {
    $src = $rowSource[$r]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $snapshot[$src][$c]
    }
}
